# Generate Report for Handoff
# Adds two new localization-status rows (942fa5b4-... and b0d3be6e-...)
# to the Overview, zh-cn and de-de sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$hyperlinkColor = 15570276  # BGR-encoded #6495ED ("cornflower blue") used by the workbook's existing HyperLink style

$file1 = "942fa5b4-8a8b-41cb-993c-03167ab1282a.md"
$file2 = "b0d3be6e-e3e7-4550-816c-a5ca6d59ebab.md"

$file1Url = "https://github.com/OpenLocalizationTest/oltest/blob/493b0ef4f33a5633f41957c788c5edb56104d16d/e2e/$file1"
$file2Url = "https://github.com/OpenLocalizationTest/oltest/blob/493b0ef4f33a5633f41957c788c5edb56104d16d/e2e/$file2"

$zhTarget1 = "942fa5b4-8a8b-41cb-993c-03167ab1282a.d2a9adcaee6d1c80e507967e457b0a6766e83171.zh-cn.xlf"
$zhTarget2 = "b0d3be6e-e3e7-4550-816c-a5ca6d59ebab.810e8da3d493daadbf705f153f9724feca4950e8.zh-cn.xlf"
$zhTarget1Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd9211dfc72cb3a44ae4b3bd06be07d0b3b32560/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$zhTarget1"
$zhTarget2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd9211dfc72cb3a44ae4b3bd06be07d0b3b32560/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/$zhTarget2"

$deTarget1 = "942fa5b4-8a8b-41cb-993c-03167ab1282a.d2a9adcaee6d1c80e507967e457b0a6766e83171.de-de.xlf"
$deTarget2 = "b0d3be6e-e3e7-4550-816c-a5ca6d59ebab.810e8da3d493daadbf705f153f9724feca4950e8.de-de.xlf"
$deTarget1Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c83f61747f28d7d57e4df238e5f53eff46fe38d5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$deTarget1"
$deTarget2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c83f61747f28d7d57e4df238e5f53eff46fe38d5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/$deTarget2"

$status = "Ready for handoff"
$handoffDate = "2016-25-09 10:25:54"
$zhHandoffDatetime = "2016-03-09 10:25:43"
$deHandoffDatetime = "2016-03-09 10:25:54"
$targetDate = "0001-01-01 00:00:00"
$reason = "Include"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$ws1.Range("A4").Value = $file1
$ws1.Hyperlinks.Add($ws1.Range("A4"), $file1Url, "", "", $file1) | Out-Null
$ws1.Range("A4").Font.Color = $hyperlinkColor
$ws1.Range("B4").Value = $status
$ws1.Range("C4").Value = $status
$ws1.Range("D4").Value = $handoffDate

$ws1.Range("A5").Value = $file2
$ws1.Hyperlinks.Add($ws1.Range("A5"), $file2Url, "", "", $file2) | Out-Null
$ws1.Range("A5").Font.Color = $hyperlinkColor
$ws1.Range("B5").Value = $status
$ws1.Range("C5").Value = $status
$ws1.Range("D5").Value = $handoffDate

# ---------------------------------------------------------------------------
# zh-cn sheet: Source File Name | File Extension | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------------
$ws2.Range("A4").Value = $file1
$ws2.Hyperlinks.Add($ws2.Range("A4"), $file1Url, "", "", $file1) | Out-Null
$ws2.Range("A4").Font.Color = $hyperlinkColor
$ws2.Range("B4").Value = ".md"
$ws2.Hyperlinks.Add($ws2.Range("B4"), $file1Url, "", "", ".md") | Out-Null
$ws2.Range("B4").Font.Color = $hyperlinkColor
$ws2.Range("C4").Value = $status
$ws2.Range("D4").Value = $zhTarget1
$ws2.Hyperlinks.Add($ws2.Range("D4"), $zhTarget1Url, "", "", $zhTarget1) | Out-Null
$ws2.Range("D4").Font.Color = $hyperlinkColor
$ws2.Range("E4").Value = $zhHandoffDatetime
$ws2.Range("H4").Value = $targetDate
$ws2.Range("I4").Value = $reason

$ws2.Range("A5").Value = $file2
$ws2.Hyperlinks.Add($ws2.Range("A5"), $file2Url, "", "", $file2) | Out-Null
$ws2.Range("A5").Font.Color = $hyperlinkColor
$ws2.Range("B5").Value = ".md"
$ws2.Hyperlinks.Add($ws2.Range("B5"), $file2Url, "", "", ".md") | Out-Null
$ws2.Range("B5").Font.Color = $hyperlinkColor
$ws2.Range("C5").Value = $status
$ws2.Range("D5").Value = $zhTarget2
$ws2.Hyperlinks.Add($ws2.Range("D5"), $zhTarget2Url, "", "", $zhTarget2) | Out-Null
$ws2.Range("D5").Font.Color = $hyperlinkColor
$ws2.Range("E5").Value = $zhHandoffDatetime
$ws2.Range("H5").Value = $targetDate
$ws2.Range("I5").Value = $reason

# ---------------------------------------------------------------------------
# de-de sheet: same layout as zh-cn
# ---------------------------------------------------------------------------
$ws3.Range("A4").Value = $file1
$ws3.Hyperlinks.Add($ws3.Range("A4"), $file1Url, "", "", $file1) | Out-Null
$ws3.Range("A4").Font.Color = $hyperlinkColor
$ws3.Range("B4").Value = ".md"
$ws3.Hyperlinks.Add($ws3.Range("B4"), $file1Url, "", "", ".md") | Out-Null
$ws3.Range("B4").Font.Color = $hyperlinkColor
$ws3.Range("C4").Value = $status
$ws3.Range("D4").Value = $deTarget1
$ws3.Hyperlinks.Add($ws3.Range("D4"), $deTarget1Url, "", "", $deTarget1) | Out-Null
$ws3.Range("D4").Font.Color = $hyperlinkColor
$ws3.Range("E4").Value = $deHandoffDatetime
$ws3.Range("H4").Value = $targetDate
$ws3.Range("I4").Value = $reason

$ws3.Range("A5").Value = $file2
$ws3.Hyperlinks.Add($ws3.Range("A5"), $file2Url, "", "", $file2) | Out-Null
$ws3.Range("A5").Font.Color = $hyperlinkColor
$ws3.Range("B5").Value = ".md"
$ws3.Hyperlinks.Add($ws3.Range("B5"), $file2Url, "", "", ".md") | Out-Null
$ws3.Range("B5").Font.Color = $hyperlinkColor
$ws3.Range("C5").Value = $status
$ws3.Range("D5").Value = $deTarget2
$ws3.Hyperlinks.Add($ws3.Range("D5"), $deTarget2Url, "", "", $deTarget2) | Out-Null
$ws3.Range("D5").Font.Color = $hyperlinkColor
$ws3.Range("E5").Value = $deHandoffDatetime
$ws3.Range("H5").Value = $targetDate
$ws3.Range("I5").Value = $reason

Write-Host "Report for Handoff rows added."
